$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unhide all previously-hidden rows first (11,12 and 41-45) before writing values ---
$ws.Range("A11:A12").EntireRow.Hidden = $false
$ws.Range("A41:A45").EntireRow.Hidden = $false

# --- Week / period header updates ---
$ws.Range("C9").Value = 5

# --- Week-start date for first block; dependent formulas (B12:B45) recalc automatically ---
$ws.Range("B11").Value = 42856

# --- Week-number labels (A11, A18, A25, A32, A39): +13 weeks ---
$ws.Range("A11").Value = 18
$ws.Range("A18").Value = 19
$ws.Range("A25").Value = 20
$ws.Range("A32").Value = 21
$ws.Range("A39").Value = 22

# --- Daily presence (D) and activity description (F) cells ---
$ws.Range("D11").Value = ""
$ws.Range("F11").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("F12").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("F13").Value = ""
$ws.Range("D14").Value = 1
$ws.Range("F14").Value = "Weryfikacja protokołów z opisania stanu nieruchomości;"
$ws.Range("D15").Value = 1
$ws.Range("F15").Value = "Weryfikacja protokołów z opisania stanu nieruchomości;"
$ws.Range("D16").Value = ""
$ws.Range("F16").Value = ""
$ws.Range("D17").Value = ""
$ws.Range("F17").Value = ""
$ws.Range("D18").Value = 1
$ws.Range("F18").Value = "Weryfikacja protokołów z opisania stanu nieruchomości;"
$ws.Range("D19").Value = 1
$ws.Range("F19").Value = "Weryfikacja protokołów z opisania stanu nieruchomości;`n; Przygotowanie do spotkania z właścicielem dz. 1020 obr. Zębowo;"
$ws.Range("D20").Value = 1
$ws.Range("F20").Value = "Weryfikacja protokołów z opisania stanu nieruchomości;"
$ws.Range("D21").Value = 1
$ws.Range("F21").Value = "Weryfikacja protokołów z opisania stanu nieruchomości;"
$ws.Range("D22").Value = 1
$ws.Range("F22").Value = "Weryfikacja protokołów z opisania stanu nieruchomości;"
$ws.Range("D23").Value = ""
$ws.Range("F23").Value = ""
$ws.Range("D24").Value = ""
$ws.Range("F24").Value = ""
$ws.Range("D25").Value = 1
$ws.Range("F25").Value = "Weryfikacja protokołów z opisania stanu nieruchomości;"
$ws.Range("D26").Value = 1
$ws.Range("F26").Value = "Weryfikacja protokołów z opisania stanu nieruchomości;`n; Wsparcie WRB w przygotowaniu protokołów przekazania drewna;"
$ws.Range("D27").Value = 1
$ws.Range("F27").Value = "Weryfikacja protokołów przekazania drewna;"
$ws.Range("D28").Value = 1
$ws.Range("F28").Value = "Weryfikacja protokołów przekazania drewna;"
$ws.Range("D29").Value = 1
$ws.Range("F29").Value = "Weryfikacja protokołów z opisania stanu nieruchomości;"
$ws.Range("D30").Value = ""
$ws.Range("F30").Value = ""
$ws.Range("D31").Value = ""
$ws.Range("F31").Value = ""
$ws.Range("D32").Value = 1
$ws.Range("F32").Value = "Weryfikacja protokołów z opisania stanu nieruchomości;"
$ws.Range("D33").Value = 1
$ws.Range("F33").Value = "Weryfikacja protokołów z opisania stanu nieruchomości;"
$ws.Range("D34").Value = 1
$ws.Range("F34").Value = "Weryfikacja protokołów z opisania stanu nieruchomości;`n; Narada koordynacyjna;"
$ws.Range("D35").Value = 1
$ws.Range("F35").Value = "Weryfikacja protokołów z opisania stanu nieruchomości;"
$ws.Range("D36").Value = 1
$ws.Range("F36").Value = "Weryfikacja protokołów z opisania stanu nieruchomości;"
$ws.Range("D37").Value = ""
$ws.Range("F37").Value = ""
$ws.Range("D38").Value = ""
$ws.Range("F38").Value = ""
$ws.Range("D39").Value = 1
$ws.Range("F39").Value = "Weryfikacja protokołów z opisania stanu nieruchomości;`n; Spotkanie z p. Jabłońskich (dot. dz. 106/3 obr. Chmielinko)"
$ws.Range("D40").Value = 1
$ws.Range("F40").Value = "Weryfikacja protokołów z opisania stanu nieruchomości;`n; Weryfikacja raportu z realizacji planu komunikacji;"
$ws.Range("D41").Value = ""
$ws.Range("F41").Value = ""
$ws.Range("D42").Value = ""
$ws.Range("F42").Value = ""
$ws.Range("D43").Value = ""
$ws.Range("F43").Value = ""
$ws.Range("D44").Value = ""
$ws.Range("F44").Value = ""
$ws.Range("D45").Value = ""
$ws.Range("F45").Value = ""

# --- Total days present ---
$ws.Range("C46").Value = 19

